$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 121, shifting existing rows 121:179 down to 122:180
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new weekly data point
$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 44452
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = 100114014
$ws.Cells.Item(121, 7).Value = "Betarraga"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 500
$ws.Cells.Item(121, 11).Value = 1000
$ws.Cells.Item(121, 12).Value = 1000
$ws.Cells.Item(121, 13).Value = 1000
$ws.Cells.Item(121, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(121, 15).Value = "Región del Maule"
$ws.Cells.Item(121, 16).Value = 200
$ws.Cells.Item(121, 17).Value = 5
$ws.Cells.Item(121, 18).Value = "Hortaliza"
